# The slide had been customized with orphaned / mismatched placeholders
# (a "Zastupny symbol pro obsah 2" content placeholder with idx=10, a
# "Nadpis 1" title, and a "Text Placeholder 2" with idx=15) that no longer
# correspond to any placeholder defined on the slide's layout
# ("Title and Content"). Fix the slide-cloning bug by resetting the slide
# back to fresh placeholders that match the layout (Title + Content
# idx=1), as PowerPoint does when a slide's placeholders are reconciled
# with its layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Resolve the slide's layout from the slide master BEFORE deleting the
# slide, via the master's CustomLayouts collection (not via $s.CustomLayout,
# which goes stale once $s is deleted).
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(2)

$s.Delete()
$newSlide = $p.Slides.AddSlide(1, $layout)

# Title placeholder -> "Title 1" (id 2, type="title")
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Click to edit Master title style"

# Content placeholder -> "Content Placeholder 2" (id 3, idx="1")
$content = $newSlide.Shapes.Item(2)
$tr = $content.TextFrame.TextRange
$tr.Text = "Edit Master text styles"
[void]$tr.InsertAfter("`rSecond level")
[void]$tr.InsertAfter("`rThird level")
[void]$tr.InsertAfter("`rFourth level")
[void]$tr.InsertAfter("`rFifth level")
for ($i = 1; $i -le 5; $i++) {
    $tr.Paragraphs($i, 1).IndentLevel = $i
}
